$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet "Recommandations": update rows 2-44 (all columns A-G) ---
# Row 2
$ws1.Range("A2").Value = "CFAO MOTORS CI"
$ws1.Range("B2").Value = 0
$ws1.Range("C2").Value = 3
$ws1.Range("D2").Value = 2760
$ws1.Range("E2").Value = 900
$ws1.Range("F2").Value = "🟡 Observer"
$ws1.Range("G2").Value = "➖ Neutre"

# Row 3
$ws1.Range("A3").Value = "UNIWAX CI"
$ws1.Range("B3").Value = 0
$ws1.Range("C3").Value = 3
$ws1.Range("D3").Value = 2515
$ws1.Range("E3").Value = 900
$ws1.Range("F3").Value = "🟡 Observer"
$ws1.Range("G3").Value = "➖ Neutre"

# Row 4
$ws1.Range("A4").Value = "BRVM - SERVICES PUBLICS"
$ws1.Range("B4").Value = 0
$ws1.Range("C4").Value = 6
$ws1.Range("D4").Value = 2485.13
$ws1.Range("E4").Value = 106.04
$ws1.Range("F4").Value = "🟡 Observer"
$ws1.Range("G4").Value = "➖ Neutre"

# Row 5
$ws1.Range("A5").Value = "SETAO CI"
$ws1.Range("B5").Value = 0
$ws1.Range("C5").Value = 3
$ws1.Range("D5").Value = 2115
$ws1.Range("E5").Value = 710
$ws1.Range("F5").Value = "🟡 Observer"
$ws1.Range("G5").Value = "➖ Neutre"

# Row 6
$ws1.Range("A6").Value = "NEI-CEDA CI"
$ws1.Range("B6").Value = 0
$ws1.Range("C6").Value = 3
$ws1.Range("D6").Value = 2025
$ws1.Range("E6").Value = 660
$ws1.Range("F6").Value = "🟡 Observer"
$ws1.Range("G6").Value = "➖ Neutre"

# Row 7
$ws1.Range("A7").Value = "BRVM - AUTRES SECTEURS"
$ws1.Range("B7").Value = 0
$ws1.Range("C7").Value = 3
$ws1.Range("D7").Value = 1851.24
$ws1.Range("E7").Value = 616.51
$ws1.Range("F7").Value = "🟡 Observer"
$ws1.Range("G7").Value = "➖ Neutre"

# Row 8
$ws1.Range("A8").Value = "AIR LIQUIDE CI"
$ws1.Range("B8").Value = 0
$ws1.Range("C8").Value = 3
$ws1.Range("D8").Value = 1780
$ws1.Range("E8").Value = 590
$ws1.Range("F8").Value = "🟡 Observer"
$ws1.Range("G8").Value = "➖ Neutre"

# Row 9
$ws1.Range("A9").Value = "BRVM - DISTRIBUTION"
$ws1.Range("B9").Value = 0
$ws1.Range("C9").Value = 3
$ws1.Range("D9").Value = 1187.17
$ws1.Range("E9").Value = 395.36
$ws1.Range("F9").Value = "🟡 Observer"
$ws1.Range("G9").Value = "➖ Neutre"

# Row 10
$ws1.Range("A10").Value = "BRVM - TRANSPORT"
$ws1.Range("B10").Value = 0
$ws1.Range("C10").Value = 3
$ws1.Range("D10").Value = 1045.18
$ws1.Range("E10").Value = 347.58
$ws1.Range("F10").Value = "🟡 Observer"
$ws1.Range("G10").Value = "➖ Neutre"

# Row 11
$ws1.Range("A11").Value = "BRVM - AGRICULTURE"
$ws1.Range("B11").Value = 0
$ws1.Range("C11").Value = 3
$ws1.Range("D11").Value = 990.88
$ws1.Range("E11").Value = 329.53
$ws1.Range("F11").Value = "🟡 Observer"
$ws1.Range("G11").Value = "➖ Neutre"

# Row 12
$ws1.Range("A12").Value = "Marché"
$ws1.Range("B12").Value = 0
$ws1.Range("C12").Value = 1
$ws1.Range("D12").Value = 805
$ws1.Range("E12").Value = 805
$ws1.Range("F12").Value = "🟡 Observer"
$ws1.Range("G12").Value = "➖ Neutre"

# Row 13
$ws1.Range("A13").Value = "BRVM - INDUSTRIELS"
$ws1.Range("B13").Value = 0
$ws1.Range("C13").Value = 3
$ws1.Range("D13").Value = 415.01
$ws1.Range("E13").Value = 138.45
$ws1.Range("F13").Value = "🟡 Observer"
$ws1.Range("G13").Value = "➖ Neutre"

# Row 14
$ws1.Range("A14").Value = "BRVM-PRESTIGE"
$ws1.Range("B14").Value = 0
$ws1.Range("C14").Value = 3
$ws1.Range("D14").Value = 402.04
$ws1.Range("E14").Value = 134.19
$ws1.Range("F14").Value = "🟡 Observer"
$ws1.Range("G14").Value = "➖ Neutre"

# Row 15
$ws1.Range("A15").Value = "BRVM - FINANCES"
$ws1.Range("B15").Value = 0
$ws1.Range("C15").Value = 3
$ws1.Range("D15").Value = 386.19
$ws1.Range("E15").Value = 129.45
$ws1.Range("F15").Value = "🟡 Observer"
$ws1.Range("G15").Value = "➖ Neutre"

# Row 16
$ws1.Range("A16").Value = "BRVM - SERVICES FINANCIERS"
$ws1.Range("B16").Value = 0
$ws1.Range("C16").Value = 3
$ws1.Range("D16").Value = 379.54
$ws1.Range("E16").Value = 127.23
$ws1.Range("F16").Value = "🟡 Observer"
$ws1.Range("G16").Value = "➖ Neutre"

# Row 17
$ws1.Range("A17").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Range("B17").Value = 0
$ws1.Range("C17").Value = 3
$ws1.Range("D17").Value = 369
$ws1.Range("E17").Value = 123.61
$ws1.Range("F17").Value = "🟡 Observer"
$ws1.Range("G17").Value = "➖ Neutre"

# Row 18
$ws1.Range("A18").Value = "BRVM - ENERGIE"
$ws1.Range("B18").Value = 0
$ws1.Range("C18").Value = 3
$ws1.Range("D18").Value = 323.71
$ws1.Range("E18").Value = 107.83
$ws1.Range("F18").Value = "🟡 Observer"
$ws1.Range("G18").Value = "➖ Neutre"

# Row 19
$ws1.Range("A19").Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Range("B19").Value = 0
$ws1.Range("C19").Value = 3
$ws1.Range("D19").Value = 283.12
$ws1.Range("E19").Value = 94.26
$ws1.Range("F19").Value = "🟡 Observer"
$ws1.Range("G19").Value = "➖ Neutre"

# Row 20
$ws1.Range("A20").Value = "BRVM - INDUSTRIE                (**)"
$ws1.Range("B20").Value = 0
$ws1.Range("C20").Value = 1
$ws1.Range("D20").Value = 244.24
$ws1.Range("E20").Value = 244.24
$ws1.Range("F20").Value = "🟡 Observer"
$ws1.Range("G20").Value = "➖ Neutre"

# Row 21
$ws1.Range("A21").Value = "BRVM - CONSOMMATION DE BASE             (**)"
$ws1.Range("B21").Value = 0
$ws1.Range("C21").Value = 1
$ws1.Range("D21").Value = 204.59
$ws1.Range("E21").Value = 204.59
$ws1.Range("F21").Value = "🟡 Observer"
$ws1.Range("G21").Value = "➖ Neutre"

# Row 22
$ws1.Range("A22").Value = "BRVM-PRINCIPAL                  (**)"
$ws1.Range("B22").Value = 0
$ws1.Range("C22").Value = 1
$ws1.Range("D22").Value = 192.37
$ws1.Range("E22").Value = 192.37
$ws1.Range("F22").Value = "🟡 Observer"
$ws1.Range("G22").Value = "➖ Neutre"

# Row 23
$ws1.Range("A23").Value = "ORAGROUP TOGO (ORGT)"
$ws1.Range("B23").Value = 3
$ws1.Range("C23").Value = 0
$ws1.Range("D23").Value = 13.05
$ws1.Range("E23").Value = 3.27
$ws1.Range("F23").Value = "🟢 Achat"
$ws1.Range("G23").Value = "✅ Renforcer"

# Row 24
$ws1.Range("A24").Value = "SMB CI (SMBC)"
$ws1.Range("B24").Value = 1
$ws1.Range("C24").Value = 0
$ws1.Range("D24").Value = 7.07
$ws1.Range("E24").Value = 7.07
$ws1.Range("F24").Value = "🟡 Observer"
$ws1.Range("G24").Value = "➖ Neutre"

# Row 25
$ws1.Range("A25").Value = "UNIWAX CI (UNXC)"
$ws1.Range("B25").Value = 2
$ws1.Range("C25").Value = 1
$ws1.Range("D25").Value = 6.92
$ws1.Range("E25").Value = 6.94
$ws1.Range("F25").Value = "🟡 Observer"
$ws1.Range("G25").Value = "👀 À surveiller"

# Row 26
$ws1.Range("A26").Value = "SOGB CI (SOGC)"
$ws1.Range("B26").Value = 1
$ws1.Range("C26").Value = 0
$ws1.Range("D26").Value = 4.53
$ws1.Range("E26").Value = 4.53
$ws1.Range("F26").Value = "🟡 Observer"
$ws1.Range("G26").Value = "➖ Neutre"

# Row 27
$ws1.Range("A27").Value = "BANK OF AFRICA CI (BOAC)"
$ws1.Range("B27").Value = 1
$ws1.Range("C27").Value = 0
$ws1.Range("D27").Value = 4.27
$ws1.Range("E27").Value = 4.27
$ws1.Range("F27").Value = "🟡 Observer"
$ws1.Range("G27").Value = "➖ Neutre"

# Row 28
$ws1.Range("A28").Value = "AIR LIQUIDE CI (SIVC)"
$ws1.Range("B28").Value = 1
$ws1.Range("C28").Value = 0
$ws1.Range("D28").Value = 4.17
$ws1.Range("E28").Value = 4.17
$ws1.Range("F28").Value = "🟡 Observer"
$ws1.Range("G28").Value = "➖ Neutre"

# Row 29
$ws1.Range("A29").Value = "BANK OF AFRICA SENEGAL (BOAS)"
$ws1.Range("B29").Value = 1
$ws1.Range("C29").Value = 0
$ws1.Range("D29").Value = 3.98
$ws1.Range("E29").Value = 3.98
$ws1.Range("F29").Value = "🟡 Observer"
$ws1.Range("G29").Value = "➖ Neutre"

# Row 30
$ws1.Range("A30").Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws1.Range("B30").Value = 1
$ws1.Range("C30").Value = 0
$ws1.Range("D30").Value = 3.23
$ws1.Range("E30").Value = 3.23
$ws1.Range("F30").Value = "🟡 Observer"
$ws1.Range("G30").Value = "➖ Neutre"

# Row 31
$ws1.Range("A31").Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Range("B31").Value = 1
$ws1.Range("C31").Value = 0
$ws1.Range("D31").Value = 3.14
$ws1.Range("E31").Value = 3.14
$ws1.Range("F31").Value = "🟡 Observer"
$ws1.Range("G31").Value = "➖ Neutre"

# Row 32
$ws1.Range("A32").Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws1.Range("B32").Value = 1
$ws1.Range("C32").Value = 0
$ws1.Range("D32").Value = 3.01
$ws1.Range("E32").Value = 3.01
$ws1.Range("F32").Value = "🟡 Observer"
$ws1.Range("G32").Value = "➖ Neutre"

# Row 33
$ws1.Range("A33").Value = "SICOR CI (SICC)"
$ws1.Range("B33").Value = 1
$ws1.Range("C33").Value = 1
$ws1.Range("D33").Value = 0.01
$ws1.Range("E33").Value = -7.42
$ws1.Range("F33").Value = "🟡 Observer"
$ws1.Range("G33").Value = "👀 À surveiller"

# Row 34
$ws1.Range("A34").Value = "TOTAL"
$ws1.Range("B34").Value = 0
$ws1.Range("C34").Value = 3
$ws1.Range("D34").Value = 0
$ws1.Range("E34").Value = 0
$ws1.Range("F34").Value = "🟡 Observer"
$ws1.Range("G34").Value = "➖ Neutre"

# Row 35
$ws1.Range("A35").Value = "NESTLE CI (NTLC)"
$ws1.Range("B35").Value = 0
$ws1.Range("C35").Value = 1
$ws1.Range("D35").Value = -1.41
$ws1.Range("E35").Value = -1.41
$ws1.Range("F35").Value = "🟡 Observer"
$ws1.Range("G35").Value = "➖ Neutre"

# Row 36
$ws1.Range("A36").Value = "SERVAIR ABIDJAN CI (ABJC)"
$ws1.Range("B36").Value = 0
$ws1.Range("C36").Value = 1
$ws1.Range("D36").Value = -1.69
$ws1.Range("E36").Value = -1.69
$ws1.Range("F36").Value = "🟡 Observer"
$ws1.Range("G36").Value = "➖ Neutre"

# Row 37
$ws1.Range("A37").Value = "NEI-CEDA CI (NEIC)"
$ws1.Range("B37").Value = 0
$ws1.Range("C37").Value = 1
$ws1.Range("D37").Value = -2.22
$ws1.Range("E37").Value = -2.22
$ws1.Range("F37").Value = "🟡 Observer"
$ws1.Range("G37").Value = "➖ Neutre"

# Row 38
$ws1.Range("A38").Value = "TOTALENERGIES MARKETING SN (TTLS)"
$ws1.Range("B38").Value = 0
$ws1.Range("C38").Value = 1
$ws1.Range("D38").Value = -2.34
$ws1.Range("E38").Value = -2.34
$ws1.Range("F38").Value = "🟡 Observer"
$ws1.Range("G38").Value = "➖ Neutre"

# Row 39
$ws1.Range("A39").Value = "BICI CI (BICC)"
$ws1.Range("B39").Value = 1
$ws1.Range("C39").Value = 1
$ws1.Range("D39").Value = -4.03
$ws1.Range("E39").Value = -7.47
$ws1.Range("F39").Value = "🟡 Observer"
$ws1.Range("G39").Value = "👀 À surveiller"

# Row 40
$ws1.Range("A40").Value = "ONATEL BF (ONTBF)"
$ws1.Range("B40").Value = 0
$ws1.Range("C40").Value = 1
$ws1.Range("D40").Value = -6.04
$ws1.Range("E40").Value = -6.04
$ws1.Range("F40").Value = "🟡 Observer"
$ws1.Range("G40").Value = "➖ Neutre"

# Row 41
$ws1.Range("A41").Value = "CFAO MOTORS CI (CFAC)"
$ws1.Range("B41").Value = 0
$ws1.Range("C41").Value = 1
$ws1.Range("D41").Value = -6.22
$ws1.Range("E41").Value = -6.22
$ws1.Range("F41").Value = "🟡 Observer"
$ws1.Range("G41").Value = "➖ Neutre"

# Row 42
$ws1.Range("A42").Value = "BERNABE CI (BNBC)"
$ws1.Range("B42").Value = 0
$ws1.Range("C42").Value = 1
$ws1.Range("D42").Value = -7.39
$ws1.Range("E42").Value = -7.39
$ws1.Range("F42").Value = "🟡 Observer"
$ws1.Range("G42").Value = "➖ Neutre"

# Row 43
$ws1.Range("A43").Value = "SAFCA CI (SAFC)"
$ws1.Range("B43").Value = 0
$ws1.Range("C43").Value = 2
$ws1.Range("D43").Value = -7.4
$ws1.Range("E43").Value = -4.17
$ws1.Range("F43").Value = "🟡 Observer"
$ws1.Range("G43").Value = "➖ Neutre"

# Row 44
$ws1.Range("A44").Value = "UNILEVER CI (UNLC)"
$ws1.Range("B44").Value = 0
$ws1.Range("C44").Value = 3
$ws1.Range("D44").Value = -21.66
$ws1.Range("E44").Value = -7.49
$ws1.Range("F44").Value = "🔴 Vente"
$ws1.Range("G44").Value = "⚠️ Risque de décrochage"

# Remove now-obsolete trailing rows 45 and 46 (data now ends at row 44)
$ws1.Rows.Item(45).Delete()
$ws1.Rows.Item(45).Delete()

# --- Sheet "Top_YTD": update column B values ---
$ws2.Range("B2").Value = 479802.99
$ws2.Range("B3").Value = 105900
$ws2.Range("B4").Value = 82255
$ws2.Range("B6").Value = 46434.8
$ws2.Range("B7").Value = 36772.29
$ws2.Range("B8").Value = 33227
$ws2.Range("B9").Value = 12081.98
$ws2.Range("B10").Value = 8915.22
$ws2.Range("B11").Value = 7866.94
